# Generate Report for Handoff
#
# The localization status report moved from "In Translation" to
# "Ready for handoff": the status cells and their corresponding
# "last generated" timestamps are refreshed on all three sheets, and
# the (now wider) status columns are resized to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Overview sheet --------------------------------------------------
# E2 / F2 hold the per-language status, G2 the latest HO xliff generate date
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-30 17:09:15"

# --- zh-cn sheet -------------------------------------------------------
# C2 holds the status, H2 the latest handoff datetime
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-30 17:09:09"

# --- de-de sheet -------------------------------------------------------
# C2 holds the status, H2 the latest handoff datetime
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-30 17:09:15"

# --- Widen the status columns to fit "Ready for handoff" --------------
# (was ~13.41 characters, grows to ~17.22 characters). The COM layer only
# keeps 1/6-character precision, so 16.33 is the input that lands on the
# nearest representable width to the target 17.2159881591797.
$overview.Columns.Item(5).ColumnWidth = 16.33   # column E
$overview.Columns.Item(6).ColumnWidth = 16.33   # column F
$zhcn.Columns.Item(3).ColumnWidth = 16.33        # column C
$dede.Columns.Item(3).ColumnWidth = 16.33        # column C
